$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 267.8
$ws.Range("I4").Value = 274.75
$ws.Range("J4").Value = 240
$ws.Range("K4").Value = 274.75
$ws.Range("L4").Value = 240
$ws.Range("M4").Value = -160.75
$ws.Range("N4").Value = -468
$ws.Range("H96").Value = 1101.8572
$ws.Range("I96").Value = 1101.8572
$ws.Range("K96").Value = 3305.5716
$ws.Range("M96").Value = -1932.5716
$ws.Range("H100").Value = 303.5
$ws.Range("I100").Value = 303.5
$ws.Range("K100").Value = 303.5
$ws.Range("M100").Value = 237.5
$ws.Range("H127").Value = 5404
$ws.Range("I127").Value = 5555.6665
$ws.Range("K127").Value = 16666.9995
$ws.Range("M127").Value = -11706.9995
$ws.Range("H138").Value = 3602.3845
$ws.Range("J138").Value = 4183.857
$ws.Range("L138").Value = 12551.571
$ws.Range("N138").Value = -22831.571
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4655.3706
$ws.Range("I32").Value = 3229.7234
$ws.Range("K32").Value = 3229.7234
$ws.Range("M32").Value = -2942.7234
$ws.Range("H74").Value = 536.46875
$ws.Range("I74").Value = 536.46875
$ws.Range("K74").Value = 536.46875
$ws.Range("M74").Value = 337.53125
$ws.Range("H77").Value = 536.46875
$ws.Range("I77").Value = 536.46875
$ws.Range("K77").Value = 2682.34375
$ws.Range("M77").Value = 1685.65625
$ws.Range("H97").Value = 639.05
$ws.Range("I97").Value = 639.64703
$ws.Range("J97").Value = 635.6667
$ws.Range("K97").Value = 639.64703
$ws.Range("L97").Value = 635.6667
$ws.Range("M97").Value = -143.64703
$ws.Range("N97").Value = -1627.6667
$ws.Range("H102").Value = 2170.2727
$ws.Range("I102").Value = 2187.3
$ws.Range("K102").Value = 2187.3
$ws.Range("M102").Value = -565.3000000000002
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864
$ws.Range("H94").Value = 646.1905
$ws.Range("I94").Value = 701.1053000000001
$ws.Range("J94").Value = 124.5
$ws.Range("K94").Value = 701.1053000000001
$ws.Range("L94").Value = 124.5
$ws.Range("M94").Value = -250.1053000000001
$ws.Range("N94").Value = -1026.5
$ws.Range("H99").Value = 1479.9688
$ws.Range("I99").Value = 1342.2858
$ws.Range("J99").Value = 2443.75
$ws.Range("K99").Value = 1342.2858
$ws.Range("L99").Value = 2443.75
$ws.Range("M99").Value = 155.7141999999999
$ws.Range("N99").Value = -5439.75
$ws.Range("H105").Value = 3141.8572
$ws.Range("I105").Value = 2826.6667
$ws.Range("K105").Value = 2826.6667
$ws.Range("M105").Value = -1079.6667
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 827.75
$ws.Range("I16").Value = 827.75
$ws.Range("K16").Value = 827.75
$ws.Range("M16").Value = -540.75
$ws.Range("H25").Value = 11667.667
$ws.Range("I25").Value = 12501
$ws.Range("J25").Value = 10001
$ws.Range("K25").Value = 12501
$ws.Range("L25").Value = 10001
$ws.Range("M25").Value = -12327
$ws.Range("N25").Value = -10349
$ws.Range("H105").Value = 3224
$ws.Range("I105").Value = 948.25
$ws.Range("K105").Value = 948.25
$ws.Range("M105").Value = 798.75
$ws.Range("H107").Value = 788.6
$ws.Range("I107").Value = 481.33334
$ws.Range("J107").Value = 1249.5
$ws.Range("K107").Value = 481.33334
$ws.Range("L107").Value = 1249.5
$ws.Range("M107").Value = 1438.66666
$ws.Range("N107").Value = -5089.5
$ws.Range("H113").Value = 827.75
$ws.Range("I113").Value = 827.75
$ws.Range("K113").Value = 827.75
$ws.Range("M113").Value = 1342.25
$ws.Range("H132").Value = 1354.0476
$ws.Range("I132").Value = 1091.3684
$ws.Range("J132").Value = 3849.5
$ws.Range("K132").Value = 3274.1052
$ws.Range("L132").Value = 11548.5
$ws.Range("M132").Value = -744.1052
$ws.Range("N132").Value = -16608.5
$ws.Range("H134").Value = 2298.1282
$ws.Range("I134").Value = 2126.2415
$ws.Range("K134").Value = 6378.7245
$ws.Range("M134").Value = -3843.7245
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1067.3334
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 600
$ws.Range("M2").Value = -487
$ws.Range("H7").Value = 33333616
$ws.Range("I7").Value = 50000050
$ws.Range("J7").Value = 751
$ws.Range("K7").Value = 150000150
$ws.Range("L7").Value = 2253
$ws.Range("M7").Value = -150000038
$ws.Range("N7").Value = -2477
$ws.Range("H34").Value = 2174.75
$ws.Range("I34").Value = 1199
$ws.Range("J34").Value = 2314.1428
$ws.Range("K34").Value = 3597
$ws.Range("L34").Value = 6942.428400000001
$ws.Range("M34").Value = -3513
$ws.Range("N34").Value = -7110.428400000001
$ws.Range("H39").Value = 2143.75
$ws.Range("I39").Value = 744
$ws.Range("J39").Value = 2983.6
$ws.Range("K39").Value = 2232
$ws.Range("L39").Value = 8950.799999999999
$ws.Range("M39").Value = -1938
$ws.Range("N39").Value = -9538.799999999999
$ws.Range("H121").Value = 398
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2298.7693
$ws.Range("I126").Value = 2409.5715
$ws.Range("K126").Value = 7228.7145
$ws.Range("M126").Value = -4758.7145
$ws.Range("H132").Value = 2692.3809
$ws.Range("I132").Value = 2042.1428
$ws.Range("J132").Value = 3992.8572
$ws.Range("K132").Value = 6126.428400000001
$ws.Range("L132").Value = 11978.5716
$ws.Range("M132").Value = -3596.428400000001
$ws.Range("N132").Value = -17038.5716
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2200
$ws.Range("I7").Value = 2200
$ws.Range("K7").Value = 2200
$ws.Range("M7").Value = -2088
$ws.Range("H16").Value = 2149.75
$ws.Range("I16").Value = 2199.6667
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 2199.6667
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -2029.6667
$ws.Range("N16").Value = -2340
$ws.Range("H22").Value = 7199.2
$ws.Range("I22").Value = 2833.3333
$ws.Range("K22").Value = 2833.3333
$ws.Range("M22").Value = -2538.3333
$ws.Range("H27").Value = 7199.2
$ws.Range("I27").Value = 2833.3333
$ws.Range("K27").Value = 2833.3333
$ws.Range("M27").Value = -2726.3333
$ws.Range("H43").Value = 8622.923000000001
$ws.Range("J43").Value = 8674.916999999999
$ws.Range("L43").Value = 8674.916999999999
$ws.Range("N43").Value = -9060.916999999999
$ws.Range("H93").Value = 1090.7142
$ws.Range("I93").Value = 943.4545000000001
$ws.Range("J93").Value = 1630.6666
$ws.Range("K93").Value = 943.4545000000001
$ws.Range("L93").Value = 1630.6666
$ws.Range("M93").Value = 304.5454999999999
$ws.Range("N93").Value = -4126.6666
$ws.Range("H126").Value = 2200
$ws.Range("I126").Value = 2200
$ws.Range("K126").Value = 6600
$ws.Range("M126").Value = -4130
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 7000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 7000
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -9746
$ws.Range("H100").Value = 3850.25
$ws.Range("I100").Value = 3799.5
$ws.Range("J100").Value = 4002.5
$ws.Range("K100").Value = 7599
$ws.Range("L100").Value = 8005
$ws.Range("M100").Value = -7058
$ws.Range("N100").Value = -9087
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

Write-Host "All updates applied"